# ADD results from server
# Update the computed results (row 2) on each year's sheet with the
# refreshed values returned by the server run.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("2025")
$ws1.Range("B2").Value = 4606.792104640001
$ws1.Range("E2").Value = 25498.12308191619
$ws1.Range("G2").Value = 6476.740570129279
$ws1.Range("I2").Value = 16841.08627624734
$ws1.Range("L2").Value = 38966.8043028528
$ws1.Range("N2").Value = 5925.595846434439
$ws1.Range("O2").Value = 5783.209061509805

$ws2 = $wb.Worksheets.Item("2030")
$ws2.Range("B2").Value = 12052.96480410071
$ws2.Range("E2").Value = 49207.63982320241
$ws2.Range("G2").Value = 6476.740570129279
$ws2.Range("I2").Value = 36930.73042578512
$ws2.Range("L2").Value = 60309.8266053096
$ws2.Range("N2").Value = 8596.626505690228
$ws2.Range("O2").Value = 6797.920363700021

$ws3 = $wb.Worksheets.Item("2035")
$ws3.Range("A2").Value = 905.2525562708205
$ws3.Range("B2").Value = 14956.7244159437
$ws3.Range("E2").Value = 62169.85694760515
$ws3.Range("G2").Value = 6476.740570129279
$ws3.Range("I2").Value = 50371.84638669913
$ws3.Range("L2").Value = 60309.8266053096
$ws3.Range("N2").Value = 12647.54260460591
$ws3.Range("O2").Value = 11316.43706107227

$ws4 = $wb.Worksheets.Item("2040")
$ws4.Range("A2").Value = 905.2525562708205
$ws4.Range("B2").Value = 14956.7244159437
$ws4.Range("E2").Value = 62169.85694760515
$ws4.Range("G2").Value = 6476.740570129279
$ws4.Range("I2").Value = 50371.84638669913
$ws4.Range("L2").Value = 60309.8266053096
$ws4.Range("N2").Value = 12647.54260460591
$ws4.Range("O2").Value = 11316.43706107227

$ws5 = $wb.Worksheets.Item("2045")
$ws5.Range("A2").Value = 905.2525562708205
$ws5.Range("B2").Value = 14956.7244159437
$ws5.Range("E2").Value = 62169.85694760515
$ws5.Range("G2").Value = 6476.740570129279
$ws5.Range("I2").Value = 50371.84638669913
$ws5.Range("L2").Value = 60309.8266053096
$ws5.Range("N2").Value = 12647.54260460591
$ws5.Range("O2").Value = 11316.43706107227

$ws6 = $wb.Worksheets.Item("2050")
$ws6.Range("A2").Value = 905.2525562708205
$ws6.Range("B2").Value = 14956.7244159437
$ws6.Range("E2").Value = 62169.85694760515
$ws6.Range("G2").Value = 6476.740570129279
$ws6.Range("I2").Value = 50371.84638669913
$ws6.Range("L2").Value = 60309.8266053096
$ws6.Range("N2").Value = 12647.54260460591
$ws6.Range("O2").Value = 11316.43706107227
